# Update "想去人数" (number of people interested) counts that changed
# between data pulls, as captured by the xml diff.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 7431
$ws1.Range("F4").Value  = 277
$ws1.Range("F5").Value  = 444
$ws1.Range("F6").Value  = 3985
$ws1.Range("F7").Value  = 319
$ws1.Range("F8").Value  = 562
$ws1.Range("F10").Value = 637
$ws1.Range("F11").Value = 121

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 3

# --- Sheet "全部类型" (All Types - combined view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 7431
$ws4.Range("F6").Value  = 277
$ws4.Range("F7").Value  = 444
$ws4.Range("F8").Value  = 3985
$ws4.Range("F9").Value  = 319
$ws4.Range("F10").Value = 562
$ws4.Range("F12").Value = 637
$ws4.Range("F13").Value = 3
$ws4.Range("F14").Value = 121
